$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Passage" column (column C), shifting Breaks and Description left.
$ws.Range("C1").EntireColumn.Delete()

# Update the selection to match the post-edit state (column C, i.e. old column D).
$ws.Range("C1:C1048576").Select()
